$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new follower row after "João Rocha" (row 31), before "Juliet" (row 32)
#    New follower: Julia Marushchenko
$ws.Rows.Item(32).Insert()
$ws.Cells.Item(32, 1).Value = "1289460729"
$ws.Cells.Item(32, 2).Value = "miss_kvitka"
$ws.Cells.Item(32, 3).Value = "Julia Marushchenko"

# After the insert above, the rows that used to be 56/57 ("Ruslan Vadimovich" and
# "Samuel Barreto") are now at 57/58. Remove those two followers - they unfollowed.
$ws.Range("A57:A58").EntireRow.Delete()

# The row that used to hold "Yulia" (id 1289460729) is now at row 70 (71 + 1 insert - 2 deletes).
# That follower's profile changed handle/display name - repurpose the row for the
# renamed/returning follower (was "Ruslan Vadimovich", now styled handle).
$ws.Cells.Item(70, 1).Value = "7239973909"
$ws.Cells.Item(70, 2).Value = "piuumee"
$ws.Cells.Item(70, 3).Value = "ʀᴜꜱʟᴀɴ ᴠᴀᴅɪᴍᴏᴠɪᴄʜ"

# Stephanie's display name lost a couple of emoji (row 62 -> now 61)
$ws.Cells.Item(61, 3).Value = "Stephanie 🌸"

# Вова Сериков's display name was updated to the full first name (row 95 -> now 94)
$ws.Cells.Item(94, 3).Value = "Владимир Сериков"

# Insert a new follower row after "Иван Безменов" (now row 107), before "Иван Пелих" (now row 108)
#    New follower: Иван Винников
$ws.Rows.Item(108).Insert()
$ws.Cells.Item(108, 1).Value = "39408861207"
$ws.Cells.Item(108, 2).Value = "ivanvinnikov_06"
$ws.Cells.Item(108, 3).Value = "Иван Винников"
